$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.555158046380211
$ws.Range("E2").Value = 2.193592700421596
$ws.Range("F2").Value = 2.825688979887564
$ws.Range("G2").Value = 3.414473754907881
$ws.Range("H2").Value = 3.94114910836752
$ws.Range("I2").Value = 4.398190732447913
$ws.Range("J2").Value = 4.784522869632266
$ws.Range("K2").Value = 5.102192146299207
$ws.Range("L2").Value = 5.354521038467411
$ws.Range("M2").Value = 5.537498246916464
$ws.Range("N2").Value = 5.656236503616492
$ws.Range("O2").Value = 5.713915123691826
$ws.Range("P2").Value = 5.711747298576652
$ws.Range("Q2").Value = 5.662266668847097
$ws.Range("R2").Value = 5.586810363776181
$ws.Range("S2").Value = 5.49931225947501
$ws.Range("T2").Value = 5.408670589460685
$ws.Range("U2").Value = 5.320372009382878
$ws.Range("V2").Value = 5.237608159387034
$ws.Range("W2").Value = 5.162044266925053
$ws.Range("X2").Value = 5.094348037817296
$ws.Range("Y2").Value = 5.034552663940652
$ws.Range("Z2").Value = 4.982304616422113
$ws.Range("AA2").Value = 4.937031233743586
$ws.Range("AB2").Value = 4.898052453675849
$ws.Range("AC2").Value = 4.864653721070636
$ws.Range("AD2").Value = 4.836132037179684
$ws.Range("AE2").Value = 4.811823577870363
$ws.Range("AF2").Value = 4.79473866154559
